$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes existing rows 3..41 down to 4..42)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the latest weekly data point
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44530
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112031
$ws.Cells.Item(3, 7).Value = "Poroto verde"
$ws.Cells.Item(3, 8).Value = "Magnum"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 40
$ws.Cells.Item(3, 11).Value = 29000
$ws.Cells.Item(3, 12).Value = 30000
$ws.Cells.Item(3, 13).Value = 29500
$ws.Cells.Item(3, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(3, 15).Value = "Región Metropolitana"
$ws.Cells.Item(3, 16).Value = 1180
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# Match the date style used by the other rows in column D (numFmtId 165 - custom datetime)
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
